# Populate Sheet1 with the uploaded content.
# Shared-string creation order matters (it reproduces the sharedStrings.xml
# table order from the target file: 0="点点滴滴", 1="的", 2=" 的"), so the
# writes below are ordered to match that.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("G10").Value = "点点滴滴"
$ws.Range("A10").Value = "的"

# Row 15
$ws.Range("O15").Value = "的"

# Row 16
$ws.Range("F16").Value = "的"

# Row 21
$ws.Range("G21").Value = "的"

# Row 22 (note the leading space in the text)
$ws.Range("C22").Value = " 的"

# Row 24
$ws.Range("M24").Value = "的"

# Row 37
$ws.Range("F37").Value = "的"

# Row 41
$ws.Range("D41").Value = " 的"
$ws.Range("L41").Value = " 的"

# Match the saved selection/active cell from the source workbook.
$ws.Range("L41").Select()
